$wb = $excel.ActiveWorkbook

# "createUser" sheet: bump the increment used to build the test user's
# username / email (1039 -> 1040). The formulas in B2/F2 recalc automatically.
$createUser = $wb.Worksheets.Item("createUser")
$createUser.Range("A2").Value = 1040

# "addListItem" sheet: the list item text mirrors that increment pattern
# (UsertwelveF -> UsertwelveG).
$addListItem = $wb.Worksheets.Item("addListItem")
$addListItem.Range("A2").Value = "UsertwelveG"

# The active/selected tab moves from "createUser" to "addListItem".
$addListItem.Activate()
